$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.739.37"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "1.898.55"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4931"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2956"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06791"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.29%  "
$ws.Range("D10").Value = "1.896.21"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07261"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("E13").Value = "  +6.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6812"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.051"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").Value = "30.729.57"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008029"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9990"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").Value = "2.138.90"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9992"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.819"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "192.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +34.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.146"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.427"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.914"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.408"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.338"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09109"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.021"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05241"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7477"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.761"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01847"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.143"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9396"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4427"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.776"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.634"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1353"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05871"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.790"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3955"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.426"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.08%  "
$ws.Range("E51").Value = "  +3.03%  "
